$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels (bold, matching the other stat headers like D6/E6, D9/E9, ...)
$ws.Range("D18").Value = "Mean increase"
$ws.Range("F18").Value = "Median increase"
$ws.Range("D18").Font.Bold = $true
$ws.Range("F18").Font.Bold = $true

# Mean / Median increase formulas
$ws.Range("D19").Formula = "=((E3 / 95.321842) * 100) - 100"
$ws.Range("D19").ClearFormats()
$ws.Range("F19").Formula = "=((E10 / 95.22216) * 100) - 100"
$ws.Range("F19").ClearFormats()

# Re-apply bold to the headers (ClearFormats above only touched D19/F19, not D18/F18)
$ws.Range("D18").Font.Bold = $true
$ws.Range("F18").Font.Bold = $true

# Update the active selection to match the authored state
$ws.Range("D21").Select()
